$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 30 new response rows (15..44) logged on 27/12/2025, appended below the
# existing history (rows 16..45). Mirrors the source sheet, where the
# "questao_id" column (C) sometimes stores the id as text instead of a
# number - NumberFormat "@" + Style "Normal" forces text storage without
# leaving a lingering text-format style on the cell.

$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = "27/12/2025 05:40"
$ws.Cells.Item(16, 3).Value = 937
$ws.Cells.Item(16, 4).Value = "Inglês"
$ws.Cells.Item(16, 5).Value = "Pronouns"
$ws.Cells.Item(16, 6).Value = 1

$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "27/12/2025 05:41"
$ws.Cells.Item(17, 3).Value = 936
$ws.Cells.Item(17, 4).Value = "Inglês"
$ws.Cells.Item(17, 5).Value = "Interpretação de Texto"
$ws.Cells.Item(17, 6).Value = 1

$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "27/12/2025 05:42"
$ws.Cells.Item(18, 3).Value = 935
$ws.Cells.Item(18, 4).Value = "Inglês"
$ws.Cells.Item(18, 5).Value = "Semantic"
$ws.Cells.Item(18, 6).Value = 0

$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "27/12/2025 05:45"
$ws.Cells.Item(19, 3).Value = 934
$ws.Cells.Item(19, 4).Value = "Inglês"
$ws.Cells.Item(19, 5).Value = "Semantic"
$ws.Cells.Item(19, 6).Value = 1

$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "27/12/2025 05:47"
$ws.Cells.Item(20, 3).Value = 933
$ws.Cells.Item(20, 4).Value = "Inglês"
$ws.Cells.Item(20, 5).Value = "Interpretação de Texto"
$ws.Cells.Item(20, 6).Value = 0

$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "27/12/2025 05:50"
$ws.Cells.Item(21, 3).NumberFormat = "@"
$ws.Cells.Item(21, 3).Value = "980"
$ws.Cells.Item(21, 3).Style = "Normal"
$ws.Cells.Item(21, 4).Value = "Inglês"
$ws.Cells.Item(21, 5).Value = "Interpretação de Texto"
$ws.Cells.Item(21, 6).Value = 1

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "27/12/2025 05:52"
$ws.Cells.Item(22, 3).Value = 979
$ws.Cells.Item(22, 4).Value = "Inglês"
$ws.Cells.Item(22, 5).Value = "Interpretação de Texto"
$ws.Cells.Item(22, 6).Value = 1

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "27/12/2025 06:01"
$ws.Cells.Item(23, 3).Value = 978
$ws.Cells.Item(23, 4).Value = "Inglês"
$ws.Cells.Item(23, 5).Value = "Interpretação de Texto"
$ws.Cells.Item(23, 6).Value = 0

$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "27/12/2025 06:01"
$ws.Cells.Item(24, 3).Value = 977
$ws.Cells.Item(24, 4).Value = "Inglês"
$ws.Cells.Item(24, 5).Value = "Semantic"
$ws.Cells.Item(24, 6).Value = 1

$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = "27/12/2025 06:02"
$ws.Cells.Item(25, 3).Value = 976
$ws.Cells.Item(25, 4).Value = "Inglês"
$ws.Cells.Item(25, 5).Value = "Semantic"
$ws.Cells.Item(25, 6).Value = 0

$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = "27/12/2025 06:06"
$ws.Cells.Item(26, 3).Value = 296
$ws.Cells.Item(26, 4).Value = "Português"
$ws.Cells.Item(26, 5).Value = "Coesão"
$ws.Cells.Item(26, 6).Value = 1

$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = "27/12/2025 06:10"
$ws.Cells.Item(27, 3).NumberFormat = "@"
$ws.Cells.Item(27, 3).Value = "334"
$ws.Cells.Item(27, 3).Style = "Normal"
$ws.Cells.Item(27, 4).Value = "Português"
$ws.Cells.Item(27, 5).Value = "Compreensão E Interpretação"
$ws.Cells.Item(27, 6).Value = 1

$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = "27/12/2025 06:11"
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "330"
$ws.Cells.Item(28, 3).Style = "Normal"
$ws.Cells.Item(28, 4).Value = "Português"
$ws.Cells.Item(28, 5).Value = "Compreensão E Interpretação"
$ws.Cells.Item(28, 6).Value = 1

$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = "27/12/2025 06:12"
$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = "329"
$ws.Cells.Item(29, 3).Style = "Normal"
$ws.Cells.Item(29, 4).Value = "Português"
$ws.Cells.Item(29, 5).Value = "Compreensão E Interpretação"
$ws.Cells.Item(29, 6).Value = 1

$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = "27/12/2025 06:16"
$ws.Cells.Item(30, 3).NumberFormat = "@"
$ws.Cells.Item(30, 3).Value = "301"
$ws.Cells.Item(30, 3).Style = "Normal"
$ws.Cells.Item(30, 4).Value = "Português"
$ws.Cells.Item(30, 5).Value = "Reescritura"
$ws.Cells.Item(30, 6).Value = 1

$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = "27/12/2025 06:18"
$ws.Cells.Item(31, 3).Value = 85
$ws.Cells.Item(31, 4).Value = "Português"
$ws.Cells.Item(31, 5).Value = "Colocação Pronominal"
$ws.Cells.Item(31, 6).Value = 1

$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = "27/12/2025 06:27"
$ws.Cells.Item(32, 3).NumberFormat = "@"
$ws.Cells.Item(32, 3).Value = "263"
$ws.Cells.Item(32, 3).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "Português"
$ws.Cells.Item(32, 5).Value = "Sinônimo E Antônimo"
$ws.Cells.Item(32, 6).Value = 0

$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = "27/12/2025 06:29"
$ws.Cells.Item(33, 3).NumberFormat = "@"
$ws.Cells.Item(33, 3).Value = "259"
$ws.Cells.Item(33, 3).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "Português"
$ws.Cells.Item(33, 5).Value = "Denotação E Conotação"
$ws.Cells.Item(33, 6).Value = 1

$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = "27/12/2025 06:30"
$ws.Cells.Item(34, 3).NumberFormat = "@"
$ws.Cells.Item(34, 3).Value = "125"
$ws.Cells.Item(34, 3).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "Português"
$ws.Cells.Item(34, 5).Value = "Correlação Verbal"
$ws.Cells.Item(34, 6).Value = 0

$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = "27/12/2025 06:37"
$ws.Cells.Item(35, 3).NumberFormat = "@"
$ws.Cells.Item(35, 3).Value = "62"
$ws.Cells.Item(35, 3).Style = "Normal"
$ws.Cells.Item(35, 4).Value = "Português"
$ws.Cells.Item(35, 5).Value = "Conjunção"
$ws.Cells.Item(35, 6).Value = 0

$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = "27/12/2025 07:11"
$ws.Cells.Item(36, 3).NumberFormat = "@"
$ws.Cells.Item(36, 3).Value = "132"
$ws.Cells.Item(36, 3).Style = "Normal"
$ws.Cells.Item(36, 4).Value = "Português"
$ws.Cells.Item(36, 5).Value = "Funções Sintáticas"
$ws.Cells.Item(36, 6).Value = 1

$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = "27/12/2025 08:32"
$ws.Cells.Item(37, 3).NumberFormat = "@"
$ws.Cells.Item(37, 3).Value = "314"
$ws.Cells.Item(37, 3).Style = "Normal"
$ws.Cells.Item(37, 4).Value = "Português"
$ws.Cells.Item(37, 5).Value = "Compreensão E Interpretação"
$ws.Cells.Item(37, 6).Value = 0

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "27/12/2025 08:35"
$ws.Cells.Item(38, 3).NumberFormat = "@"
$ws.Cells.Item(38, 3).Value = "310"
$ws.Cells.Item(38, 3).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "Português"
$ws.Cells.Item(38, 5).Value = "Narração"
$ws.Cells.Item(38, 6).Value = 1

$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = "27/12/2025 08:37"
$ws.Cells.Item(39, 3).NumberFormat = "@"
$ws.Cells.Item(39, 3).Value = "309"
$ws.Cells.Item(39, 3).Style = "Normal"
$ws.Cells.Item(39, 4).Value = "Português"
$ws.Cells.Item(39, 5).Value = "Narração"
$ws.Cells.Item(39, 6).Value = 1

$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = "27/12/2025 08:39"
$ws.Cells.Item(40, 3).NumberFormat = "@"
$ws.Cells.Item(40, 3).Value = "279"
$ws.Cells.Item(40, 3).Style = "Normal"
$ws.Cells.Item(40, 4).Value = "Português"
$ws.Cells.Item(40, 5).Value = "Coesão"
$ws.Cells.Item(40, 6).Value = 1

$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = "27/12/2025 08:41"
$ws.Cells.Item(41, 3).NumberFormat = "@"
$ws.Cells.Item(41, 3).Value = "278"
$ws.Cells.Item(41, 3).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "Português"
$ws.Cells.Item(41, 5).Value = "Sinônimo E Antônimo"
$ws.Cells.Item(41, 6).Value = 1

$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = "27/12/2025 08:41"
$ws.Cells.Item(42, 3).NumberFormat = "@"
$ws.Cells.Item(42, 3).Value = "262"
$ws.Cells.Item(42, 3).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "Português"
$ws.Cells.Item(42, 5).Value = "Sinônimo E Antônimo"
$ws.Cells.Item(42, 6).Value = 1

$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "27/12/2025 08:42"
$ws.Cells.Item(43, 3).NumberFormat = "@"
$ws.Cells.Item(43, 3).Value = "261"
$ws.Cells.Item(43, 3).Style = "Normal"
$ws.Cells.Item(43, 4).Value = "Português"
$ws.Cells.Item(43, 5).Value = "Sinônimo E Antônimo"
$ws.Cells.Item(43, 6).Value = 1

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "27/12/2025 08:43"
$ws.Cells.Item(44, 3).NumberFormat = "@"
$ws.Cells.Item(44, 3).Value = "233"
$ws.Cells.Item(44, 3).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "Português"
$ws.Cells.Item(44, 5).Value = "Crase"
$ws.Cells.Item(44, 6).Value = 1

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "27/12/2025 08:44"
$ws.Cells.Item(45, 3).NumberFormat = "@"
$ws.Cells.Item(45, 3).Value = "151"
$ws.Cells.Item(45, 3).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "Português"
$ws.Cells.Item(45, 5).Value = "Vírgula"
$ws.Cells.Item(45, 6).Value = 1
